# Regenerate demo files using new option "connectPoints" for chocolate
#
# The underlying data for the Posthoc_1 table was recomputed (new model run),
# so the emmCI strings, the posthoc statistics and a couple of derived
# labels change. Column A's per-row contrast id also changes (1/1/1 -> 1/2/3)
# and the D/E column widths grow slightly to fit the new, longer CI strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds the per-row contrast id as TEXT (e.g. "1", "2", "3"), not
# a number - force text formatting before assigning so Excel doesn't
# auto-convert the numeric-looking string to a numeric cell.
$ws.Range("A3:A4").NumberFormat = "@"

# --- Row 2 (A2 is already "1" - unchanged) ----------------------------
$ws.Range("D2").Value = "242.49 (201.43, 293.76)"
$ws.Range("E2").Value = "294.89 (246.87, 354.23)"
$ws.Range("F2").Value = 0.10300696953510069
$ws.Range("G2").Value = 0.20601393907020138
$ws.Range("H2").Value = 52.396985389858088
$ws.Range("I2").Value = 21.607910458832563
$ws.Range("J2").Value = 2.6674883590925078
$ws.Range("L2").Value = 534
$ws.Range("M2").Value = 0.14135482379611519
$ws.Range("N2").Value = 0.0049704675929756531
$ws.Range("O2").Value = "very small"

# --- Row 3 -----------------------------------------------------------
$ws.Range("A3").Value = "2"
$ws.Range("D3").Value = "281.47 (240.23, 331.26)"
$ws.Range("E3").Value = "311.36 (262.56, 371.12)"
$ws.Range("F3").Value = 0.40809498568816061
$ws.Range("G3").Value = 0.40809498568816061
$ws.Range("H3").Value = 29.892962547327272
$ws.Range("I3").Value = 10.620445338454356
$ws.Range("J3").Value = 0.6854242252662146
$ws.Range("L3").Value = 534
$ws.Range("M3").Value = 0.071653777706915209
$ws.Range("N3").Value = 0.0012819205353491014
$ws.Range("O3").Value = "very small"

# --- Row 4 -----------------------------------------------------------
$ws.Range("A4").Value = "3"
$ws.Range("D4").Value = "370.28 (308.01, 447.71)"
$ws.Range("E4").Value = "451.44 (380.15, 538.72)"
$ws.Range("F4").Value = 0.01463207336137261
$ws.Range("G4").Value = 0.043896220084117865
$ws.Range("H4").Value = 81.162516556102105
$ws.Range("I4").Value = 21.919505087962246
$ws.Range("J4").Value = 5.9992645534675786
$ws.Range("L4").Value = 534
$ws.Range("M4").Value = 0.21198658274034568
$ws.Range("N4").Value = 0.011109764303898541
$ws.Range("O4").Value = "small"
$ws.Range("P4").Value = "*"

# --- Column widths (D/E grow to fit the longer CI strings, H shrinks) -
$ws.Range("D1").ColumnWidth = 19.833333333333332
$ws.Range("E1").ColumnWidth = 19.833333333333332
$ws.Range("H1").ColumnWidth = 10.833333333333334
